$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 971.46155
$ws.Range("J19").Value = 1023
$ws.Range("L19").Value = 1023
$ws.Range("N19").Value = -1373
$ws.Range("H28").Value = 4957.933
$ws.Range("I28").Value = 6201.3
$ws.Range("J28").Value = 2471.2
$ws.Range("K28").Value = 6201.3
$ws.Range("L28").Value = 2471.2
$ws.Range("M28").Value = -5716.3
$ws.Range("N28").Value = -3441.2
$ws.Range("H62").Value = 111116990
$ws.Range("I62").Value = 200006380
$ws.Range("J62").Value = 5250
$ws.Range("K62").Value = 200006380
$ws.Range("L62").Value = 5250
$ws.Range("M62").Value = -200005756
$ws.Range("N62").Value = -6498
$ws.Range("H65").Value = 111116990
$ws.Range("I65").Value = 200006380
$ws.Range("J65").Value = 5250
$ws.Range("K65").Value = 1000031900
$ws.Range("L65").Value = 26250
$ws.Range("M65").Value = -1000028780
$ws.Range("N65").Value = -32490
$ws.Range("H98").Value = 42790.266
$ws.Range("I98").Value = 55744.8
$ws.Range("J98").Value = 16881.2
$ws.Range("K98").Value = 55744.8
$ws.Range("L98").Value = 16881.2
$ws.Range("M98").Value = -54246.8
$ws.Range("N98").Value = -19877.2
$ws.Range("H116").Value = 337043.56
$ws.Range("J116").Value = 3999
$ws.Range("L116").Value = 3999
$ws.Range("N116").Value = -10883
$ws.Range("H118").Value = 1647.5
$ws.Range("I118").Value = 1647.5
$ws.Range("K118").Value = 4942.5
$ws.Range("M118").Value = -3285.5
$ws.Range("H122").Value = 42790.266
$ws.Range("I122").Value = 55744.8
$ws.Range("J122").Value = 16881.2
$ws.Range("K122").Value = 167234.4
$ws.Range("L122").Value = 50643.60000000001
$ws.Range("M122").Value = -164784.4
$ws.Range("N122").Value = -55543.60000000001
$ws.Range("H123").Value = 106472.5
$ws.Range("J123").Value = 106472.5
$ws.Range("L123").Value = 106472.5
$ws.Range("N123").Value = -116272.5
$ws.Range("H138").Value = 5052.606
$ws.Range("J138").Value = 6180.4375
$ws.Range("L138").Value = 18541.3125
$ws.Range("N138").Value = -28821.3125
$ws.Range("H141").Value = 4100.65
$ws.Range("I141").Value = 2286.75
$ws.Range("J141").Value = 6821.5
$ws.Range("K141").Value = 6860.25
$ws.Range("L141").Value = 20464.5
$ws.Range("M141").Value = -1680.25
$ws.Range("N141").Value = -30824.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 3745
$ws.Range("I22").Value = 3745
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3745
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3446
$ws.Range("N22").ClearContents()
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10804
$ws.Range("H74").Value = 2707.3635
$ws.Range("I74").Value = 2478.1
$ws.Range("K74").Value = 2478.1
$ws.Range("M74").Value = -1604.1
$ws.Range("H77").Value = 2707.3635
$ws.Range("I77").Value = 2478.1
$ws.Range("K77").Value = 12390.5
$ws.Range("M77").Value = -8022.5
$ws.Range("H88").Value = 100002136
$ws.Range("J88").Value = 125002170
$ws.Range("L88").Value = 125002170
$ws.Range("N88").Value = -125002982
$ws.Range("H91").Value = 100002136
$ws.Range("J91").Value = 125002170
$ws.Range("L91").Value = 125002170
$ws.Range("N91").Value = -125004978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2966.7693
$ws.Range("I31").Value = 1561.75
$ws.Range("K31").Value = 1561.75
$ws.Range("M31").Value = -1266.75
$ws.Range("H34").Value = 2966.7693
$ws.Range("I34").Value = 1561.75
$ws.Range("K34").Value = 1561.75
$ws.Range("M34").Value = -1359.75
$ws.Range("H118").Value = 60742
$ws.Range("J118").Value = 60742
$ws.Range("L118").Value = 60742
$ws.Range("N118").Value = -64056
$ws.Range("H134").Value = 1791702.5
$ws.Range("I134").Value = 3480062.5
$ws.Range("K134").Value = 10440187.5
$ws.Range("M134").Value = -10437652.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 286994.06
$ws.Range("J5").Value = 668308.5600000001
$ws.Range("L5").Value = 2004925.68
$ws.Range("N5").Value = -2005149.68
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H68").Value = 13163874
$ws.Range("I68").Value = 2314.2727
$ws.Range("J68").Value = 18525990
$ws.Range("K68").Value = 6942.8181
$ws.Range("L68").Value = 55577970
$ws.Range("M68").Value = -6131.8181
$ws.Range("N68").Value = -55579592
$ws.Range("H71").Value = 13163874
$ws.Range("I71").Value = 2314.2727
$ws.Range("J71").Value = 18525990
$ws.Range("K71").Value = 20828.4543
$ws.Range("L71").Value = 166733910
$ws.Range("M71").Value = -16772.4543
$ws.Range("N71").Value = -166742022
$ws.Range("H80").Value = 108666.086
$ws.Range("J80").Value = 158818.12
$ws.Range("L80").Value = 476454.36
$ws.Range("N80").Value = -478326.36
$ws.Range("H83").Value = 108666.086
$ws.Range("J83").Value = 158818.12
$ws.Range("L83").Value = 1429363.08
$ws.Range("N83").Value = -1438723.08
$ws.Range("H103").Value = 50000
$ws.Range("I103").Value = 50000
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 150000
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -149121
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 1766.35
$ws.Range("J113").Value = 1791
$ws.Range("L113").Value = 5373
$ws.Range("N113").Value = -9713
$ws.Range("H121").Value = 1815672.9
$ws.Range("I121").Value = 1582617.2
$ws.Range("J121").Value = 2002117.2
$ws.Range("K121").Value = 4747851.6
$ws.Range("L121").Value = 6006351.6
$ws.Range("M121").Value = -4746541.6
$ws.Range("N121").Value = -6008971.6
$ws.Range("H134").Value = 13686.875
$ws.Range("I134").Value = 13686.875
$ws.Range("K134").Value = 41060.625
$ws.Range("M134").Value = -35990.625
$ws.Range("H135").Value = 286994.06
$ws.Range("J135").Value = 668308.5600000001
$ws.Range("L135").Value = 6014777.040000001
$ws.Range("N135").Value = -6019847.040000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11093.25
$ws.Range("I70").Value = 11093.25
$ws.Range("K70").Value = 11093.25
$ws.Range("M70").Value = -10823.25
$ws.Range("H73").Value = 11093.25
$ws.Range("I73").Value = 11093.25
$ws.Range("K73").Value = 11093.25
$ws.Range("M73").Value = -10157.25
$ws.Range("H132").Value = 4863.9614
$ws.Range("J132").Value = 8879.799999999999
$ws.Range("L132").Value = 26639.4
$ws.Range("N132").Value = -31699.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H61").Value = 2700.1875
$ws.Range("I61").Value = 2508.75
$ws.Range("K61").Value = 2508.75
$ws.Range("M61").Value = -2306.75
$ws.Range("H113").Value = 2700.1875
$ws.Range("I113").Value = 2508.75
$ws.Range("K113").Value = 2508.75
$ws.Range("M113").Value = -338.75
$ws.Range("H122").Value = 5555.3335
$ws.Range("J122").Value = 6420
$ws.Range("L122").Value = 19260
$ws.Range("N122").Value = -24160
